# 🔄 Actualización automática del mapa (2025-08-11 07:12:53)
# El reclamo "HUMAHUACA 3828" (Caso 6543, fila 82) fue resuelto/eliminado
# de la planilla; el resto de las filas se recorre una posición hacia
# arriba y la última fila (85) desaparece.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(82).Delete()
